$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from "EP-8" to "EA-8,EP-8"
$ws.Range("B9").Value = "EA-8,EP-8"
$ws.Range("C9").Value = "EA-8,EP-8"

# Remove the "Requisitos:" rows (24 and 25)
$ws.Rows.Item(24).Delete()
$ws.Rows.Item(24).Delete()
